$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Heading "Konfi" + bookmark + "guration" -> single run "Konfiguration"
#    (the stray _GoBack bookmark that split the word is removed).
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Konfiguration", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Konfiguration", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Insert "außer dem Server " right before "eingetragen werden." in the
#    sentence about entering the server's IPv4 address into the config
#    files of ALL application parts, and drop a (now current) _GoBack
#    bookmark right after the inserted text - mirroring where Word
#    leaves _GoBack after the last edit.
# ---------------------------------------------------------------------
$target = $d.Content
$target.Find.Execute("eingetragen werden.", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null

$insertionPoint = $d.Range($target.Start, $target.Start)
$d.Bookmarks.Add("_GoBack", $insertionPoint) | Out-Null

$bookmarkStart = $d.Bookmarks("_GoBack").Range.Start
$ip = $d.Range($bookmarkStart, $bookmarkStart)
$ip.Text = "außer dem Server "
